$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "60.144.67"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +5.15%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.597.12"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +7.48%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "505.77"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +3.35%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "156.60"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  -0.09%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.590"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -3.87%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.639.12"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +8.51%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "6.48"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +6.23%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.104"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +4.50%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.343"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("E13").Value = "  +0.77%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.083.68"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +8.65%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "60.296.31"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +5.45%  "
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("E17").Value = "  +4.85%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.640.33"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +8.48%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.77"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.33%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "344.33"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +6.25%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.45"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +4.40%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.16"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +3.85%  "
$ws.Range("E23").Value = "  -0.05%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "60.02"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +3.81%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.423"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +4.69%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.166"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0858"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +9.23%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.54"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +3.57%  "
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "19.45"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "156.01"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.57"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.75"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +8.55%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.03"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +7.00%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.21"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.31%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "307.13"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +8.82%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.49"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +8.21%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.848"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +3.40%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.76"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +6.79%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.834"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +28.50%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "35.32"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +3.86%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.631"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0572"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +7.81%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.100"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "19.95"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +13.48%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.993"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "4.84"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.047.20"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +7.95%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0235"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "10.27"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.56%  "
